# Add "UAT kelompok 10" header block in column I (I1:I3)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values -----------------------------------------------------------
$ws.Range("I1").Value = "KELOMPOK 10 : "
$ws.Range("I2").Value = "1. KEN ABEL VALLERON LIMANSYAH - C14210227"
$ws.Range("I3").Value = "2. JASON ENRICO SANTOSO - C14210067"

# --- Column width -------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 50.63

# --- I1: big centered title font (size 20) ------------------------------
$i1 = $ws.Range("I1")
$i1.Font.ThemeColor = 1
$i1.WrapText = $true
$i1.VerticalAlignment = -4108
$i1.HorizontalAlignment = -4108
$i1.Font.Size = 20

# --- I2: member name 1, size 18, left aligned, vertical centered --------
$i2 = $ws.Range("I2")
$i2.Font.ThemeColor = 1
$i2.WrapText = $true
$i2.VerticalAlignment = -4108
$i2.Font.Size = 18

# --- I3: member name 2, size 18, left aligned, vertical centered --------
$i3 = $ws.Range("I3")
$i3.Font.ThemeColor = 1
$i3.WrapText = $true
$i3.VerticalAlignment = -4108
$i3.Font.Size = 18
